$wb = $excel.ActiveWorkbook

# --- Sheet 1: Weekly Quantity ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$ws1.Cells.Item(15, 1).Value = 45662.99999999999
$ws1.Cells.Item(15, 2).Value = 6
$ws1.Cells.Item(16, 1).Value = 45669.99999999999
$ws1.Cells.Item(16, 2).Value = 12
$ws1.Cells.Item(17, 1).Value = 45676.99999999999
$ws1.Cells.Item(17, 2).Value = 4

$ws1.Range("A15:A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Sheet 2: Monthly Trend ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(9, 1).Value = 45688.99999999999
$ws2.Cells.Item(9, 2).Value = 22

$ws2.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Sheet 3: PO Forecast ---
$ws3 = $wb.Worksheets.Item("PO Forecast")

$ws3.Cells.Item(2, 2).Value = 156
$ws3.Cells.Item(3, 2).Value = 111
$ws3.Cells.Item(4, 2).Value = 101
$ws3.Cells.Item(5, 2).Value = 92
$ws3.Cells.Item(6, 2).Value = 88
$ws3.Cells.Item(7, 2).Value = 83
$ws3.Cells.Item(8, 2).Value = 74
$ws3.Cells.Item(9, 2).Value = 38
$ws3.Cells.Item(10, 2).Value = 33
$ws3.Cells.Item(11, 2).Value = 24
$ws3.Cells.Item(12, 2).Value = 19
$ws3.Cells.Item(13, 2).Value = 15
$ws3.Cells.Item(14, 2).Value = 10

$ws3.Cells.Item(15, 1).Value = 45662.99999999999
$ws3.Cells.Item(16, 1).Value = 45669.99999999999
$ws3.Cells.Item(17, 1).Value = 45676.99999999999
$ws3.Cells.Item(18, 1).Value = 45683.99999999999
$ws3.Cells.Item(19, 1).Value = 45690.99999999999
$ws3.Cells.Item(20, 1).Value = 45697.99999999999
$ws3.Cells.Item(21, 1).Value = 45704.99999999999
$ws3.Cells.Item(21, 2).Value = 0
$ws3.Cells.Item(22, 1).Value = 45711.99999999999
$ws3.Cells.Item(22, 2).Value = 0
$ws3.Cells.Item(23, 1).Value = 45718.99999999999
$ws3.Cells.Item(23, 2).Value = 0
$ws3.Cells.Item(24, 1).Value = 45725.99999999999
$ws3.Cells.Item(24, 2).Value = 0
$ws3.Cells.Item(25, 1).Value = 45732.99999999999
$ws3.Cells.Item(25, 2).Value = 0

$ws3.Range("A21:A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
